# feature/sound effect on pickup
# Adds a new credits row (row 3) for the pickup sound effect, and
# refreshes the license label text in row 2/3 to "Creative Commons License 0".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New source URL for the pickup sound effect (set first so this string is
# registered before the updated license text, matching the original edit order)
$ws.Range("B3").Value = "https://freesound.org/people/suntemple/sounds/253172/"

# Update the licensing text (applies to both the existing and new rows)
$ws.Range("C2").Value = "Creative Commons License 0"
$ws.Range("C3").Value = "Creative Commons License 0"

# New row's "Notes/Other" column
$ws.Range("D3").Value = "Doesn't require attribution"

# Restore the view/selection state to match where the author ended up editing
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("E3").Select()
